$d = $word.ActiveDocument

function Find-RangeFrom {
    param(
        [string]$findText,
        [int]$startAt = 0
    )
    $rng = $d.Range($startAt, $d.Content.End)
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Find failed for: $findText"
    }
    return $rng
}

function Replace-ParaXml {
    param(
        [string]$findText,
        [string]$paraXml,
        [int]$startAt = 0
    )
    $rng = Find-RangeFrom $findText $startAt
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + $paraXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
    return $rng.Start
}

$cursor = 0

# 1) "Vers" + "." -> "Vers." (version label cell)
$cursor = Replace-ParaXml "Vers." '<w:body><w:p w14:paraId="50028B67" w14:textId="77777777" w:rsidR="003F329D" w:rsidRPr="00E07CE9" w:rsidRDefault="003F329D" w:rsidP="000C05D4"><w:pPr><w:rPr><w:bCs/><w:i/><w:iCs/></w:rPr></w:pPr><w:r w:rsidRPr="00E07CE9"><w:rPr><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>Vers.</w:t></w:r></w:p></w:body>' $cursor

# 2) "Entry " + "Condition" -> "Entry Condition"
$cursor = Replace-ParaXml "Entry Condition" '<w:body><w:p w14:paraId="6A0C6ABC" w14:textId="77777777" w:rsidR="003F329D" w:rsidRPr="00E07CE9" w:rsidRDefault="003F329D" w:rsidP="000C05D4"><w:r w:rsidRPr="00E07CE9"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Entry Condition</w:t></w:r></w:p></w:body>' $cursor

# 3) "Exit " + "condition" -> "Exit condition" (first occurrence, success row)
$cursor = Replace-ParaXml "Exit condition" '<w:body><w:p w14:paraId="5BB143CF" w14:textId="77777777" w:rsidR="003F329D" w:rsidRPr="00E07CE9" w:rsidRDefault="003F329D" w:rsidP="000C05D4"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="00E07CE9"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Exit condition</w:t></w:r></w:p></w:body>' $cursor

# 4) "Exit " + "condition" -> "Exit condition" (second occurrence, failure row)
$cursor = Replace-ParaXml "Exit condition" '<w:body><w:p w14:paraId="2A132C6F" w14:textId="77777777" w:rsidR="003F329D" w:rsidRPr="00E07CE9" w:rsidRDefault="003F329D" w:rsidP="000C05D4"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="00E07CE9"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Exit condition</w:t></w:r></w:p></w:body>' $cursor

# 5) "On " + "failure" -> "On failure"
$cursor = Replace-ParaXml "On failure" '<w:body><w:p w14:paraId="3C9C62A9" w14:textId="77777777" w:rsidR="003F329D" w:rsidRPr="00E07CE9" w:rsidRDefault="003F329D" w:rsidP="000C05D4"><w:r w:rsidRPr="00E07CE9"><w:t xml:space="preserve">                       On failure</w:t></w:r></w:p></w:body>' $cursor

# 6) "Rilevanza/User " + "Priority" -> "Rilevanza/User Priority"
$cursor = Replace-ParaXml "Rilevanza/User Priority" '<w:body><w:p w14:paraId="132C7A3B" w14:textId="77777777" w:rsidR="003F329D" w:rsidRPr="00E07CE9" w:rsidRDefault="003F329D" w:rsidP="000C05D4"><w:r w:rsidRPr="00E07CE9"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Rilevanza/User Priority</w:t></w:r></w:p></w:body>' $cursor

# 7) "Generalization" + " of" -> "Generalization of"
$cursor = Replace-ParaXml "Generalization of" '<w:body><w:p w14:paraId="1CF66B30" w14:textId="77777777" w:rsidR="003F329D" w:rsidRPr="00E07CE9" w:rsidRDefault="003F329D" w:rsidP="000C05D4"><w:pPr><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="00E07CE9"><w:rPr><w:rFonts w:eastAsia="Calibri"/><w:b/><w:bCs/></w:rPr><w:t>Generalization of</w:t></w:r></w:p></w:body>' $cursor

# 8) "Il sistema effettua il " + "refresh" + " automatico..." -> single run
$cursor = Replace-ParaXml "Il sistema effettua il refresh automatico della pagina per visualizzare l’elenco aggiornato" '<w:body><w:p w14:paraId="4337E430" w14:textId="77777777" w:rsidR="003F329D" w:rsidRPr="00E07CE9" w:rsidRDefault="003F329D" w:rsidP="000C05D4"><w:r w:rsidRPr="00E07CE9"><w:t>Il sistema effettua il refresh automatico della pagina per visualizzare l&#8217;elenco aggiornato</w:t></w:r></w:p></w:body>' $cursor

# 9) "Il" + " Scenario/Flusso di eventi Alternativo:  " -> single run (gramStart/gramEnd removed)
$cursor = Replace-ParaXml "Il Scenario/Flusso di eventi Alternativo:  " '<w:body><w:p w14:paraId="5A7AF353" w14:textId="77777777" w:rsidR="003F329D" w:rsidRPr="00E07CE9" w:rsidRDefault="003F329D" w:rsidP="000C05D4"><w:r w:rsidRPr="00E07CE9"><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Il Scenario/Flusso di eventi Alternativo:  </w:t></w:r><w:r w:rsidRPr="00E07CE9"><w:t>Non &#232; presente alcuna prenotazione per quella giornata</w:t></w:r></w:p></w:body>' $cursor

# 10) Delete the extra empty table row (after the "refresh" row)
$t = $d.Tables(1)
$cell = $t.Cell(22, 1)
$cell.Row.Delete()
